$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Nagel
$ws.Range("B2").Value = "Nagel"
$ws.Range("E2").Value = 100
$ws.Range("F2").Value = 850

# Row 3 - Lucas
$ws.Range("B3").Value = "Lucas"
$ws.Range("C3").Value = 6.5
$ws.Range("F3").Value = 650

# Row 4 - Funke
$ws.Range("B4").Value = "Funke"
$ws.Range("C4").Value = 7
$ws.Range("F4").Value = 700

# Row 5 - McClure (Maint. Items)
$ws.Range("B5").Value = "McClure (Maint. Items)"
$ws.Range("C5").Value = 7
$ws.Range("F5").Value = 700

# Row 6 - McGill
$ws.Range("B6").Value = "McGill"
$ws.Range("C6").Value = 7.5
$ws.Range("F6").Value = 675

# Row 8 - SUBTOTAL
$ws.Range("C8").Value = 36.5
$ws.Range("D8").Value = "Reg: 36.5 / OT: 0"
$ws.Range("F8").Value = 3575
